# Burndown Chart Updated for October 15th
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Burn down chart")

# Update the "Sprint 2 Burndown" table (rows 13-27, column B & C)
$ws.Range("B13").Value = 58
$ws.Range("C13").Value = 58
$ws.Range("C14").Value = 57
$ws.Range("C15").Value = 55
$ws.Range("C16").Value = 54
$ws.Range("C17").Value = 52
$ws.Range("C18").Value = 50
$ws.Range("C19").Value = 47
$ws.Range("C20").Value = 45
$ws.Range("C21").Value = 42
$ws.Range("C22").Value = 39
$ws.Range("C23").Value = 34
$ws.Range("C24").Value = 29
$ws.Range("C25").Value = 25
$ws.Range("C26").Value = 22
$ws.Range("C27").Value = 18

# Update the view to match the scrolled/selected state in the diff
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("C31").Select()

$wb.Save()
